{"js": "// Auto-generated replacements for two-digit multiplication table values.\nconst replacements = [\n  [\"56\u00d791=5096\", \"26\u00d797=2522\"],\n  [\"93\u00d745=4185\", \"17\u00d760=1020\"],\n  [\"85\u00d728=2380\", \"18\u00d735=630\"],\n  [\"63\u00d796=6048\", \"30\u00d725=750\"],\n  [\"71\u00d732=2272\", \"78\u00d754=4212\"],\n  [\"45\u00d787=3915\", \"40\u00d731=1240\"],\n  [\"29\u00d749=1421\", \"80\u00d760=4800\"],\n  [\"55\u00d720=1100\", \"85\u00d733=2805\"],\n  [\"65\u00d734=2210\", \"92\u00d768=6256\"],\n  [\"15\u00d788=1320\", \"75\u00d711=825\"],\n  [\"52\u00d753=2756\", \"15\u00d729=435\"],\n  [\"26\u00d749=1274\", \"42\u00d782=3444\"],\n  [\"41\u00d722=902\", \"67\u00d742=2814\"],\n  [\"75\u00d736=2700\", \"23\u00d774=1702\"],\n  [\"56\u00d787=4872\", \"21\u00d732=672\"],\n  [\"58\u00d796=5568\", \"72\u00d793=6696\"],\n  [\"13\u00d719=247\", \"37\u00d743=1591\"],\n  [\"74\u00d790=6660\", \"56\u00d773=4088\"],\n  [\"71\u00d726=1846\", \"13\u00d755=715\"],\n  [\"37\u00d744=1628\", \"21\u00d762=1302\"],\n  [\"93\u00d791=8463\", \"60\u00d793=5580\"],\n  [\"39\u00d745=1755\", \"21\u00d789=1869\"],\n  [\"98\u00d734=3332\", \"38\u00d748=1824\"],\n  [\"70\u00d712=840\", \"53\u00d738=2014\"],\n  [\"27\u00d748=1296\", \"93\u00d760=5580\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Auto-generated replacements for two-digit multiplication table values.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"56\u00d791=5096\", \"26\u00d797=2522\"),\n  @(\"93\u00d745=4185\", \"17\u00d760=1020\"),\n  @(\"85\u00d728=2380\", \"18\u00d735=630\"),\n  @(\"63\u00d796=6048\", \"30\u00d725=750\"),\n  @(\"71\u00d732=2272\", \"78\u00d754=4212\"),\n  @(\"45\u00d787=3915\", \"40\u00d731=1240\"),\n  @(\"29\u00d749=1421\", \"80\u00d760=4800\"),\n  @(\"55\u00d720=1100\", \"85\u00d733=2805\"),\n  @(\"65\u00d734=2210\", \"92\u00d768=6256\"),\n  @(\"15\u00d788=1320\", \"75\u00d711=825\"),\n  @(\"52\u00d753=2756\", \"15\u00d729=435\"),\n  @(\"26\u00d749=1274\", \"42\u00d782=3444\"),\n  @(\"41\u00d722=902\", \"67\u00d742=2814\"),\n  @(\"75\u00d736=2700\", \"23\u00d774=1702\"),\n  @(\"56\u00d787=4872\", \"21\u00d732=672\"),\n  @(\"58\u00d796=5568\", \"72\u00d793=6696\"),\n  @(\"13\u00d719=247\", \"37\u00d743=1591\"),\n  @(\"74\u00d790=6660\", \"56\u00d773=4088\"),\n  @(\"71\u00d726=1846\", \"13\u00d755=715\"),\n  @(\"37\u00d744=1628\", \"21\u00d762=1302\"),\n  @(\"93\u00d791=8463\", \"60\u00d793=5580\"),\n  @(\"39\u00d745=1755\", \"21\u00d789=1869\"),\n  @(\"98\u00d734=3332\", \"38\u00d748=1824\"),\n  @(\"70\u00d712=840\", \"53\u00d738=2014\"),\n  @(\"27\u00d748=1296\", \"93\u00d760=5580\"),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
